$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


# Row 5
$ws.Range("I5").Value = 2.64
$ws.Range("J5").Value = 3.35
$ws.Range("V5").Value = 1.6

# Row 7
$ws.Range("H7").Value = 2.42
$ws.Range("I7").Value = 2.76

# Row 8
$ws.Range("L8").Value = 1.41

# Row 9
$ws.Range("L9").Value = 1.01

# Row 10
$ws.Range("AC10").Value = 10.5
$ws.Range("F10").Value = 2.76
$ws.Range("H10").Value = 2.58
$ws.Range("I10").Value = 2.7
$ws.Range("J10").Value = 3.65
$ws.Range("K10").Value = 4.1
$ws.Range("L10").Value = 1.29
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 3.75
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 2.08
$ws.Range("Q10").Value = 1.64
$ws.Range("R10").Value = 1.43
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.63
$ws.Range("U10").Value = 2.3
$ws.Range("V10").Value = 1.58

# Row 11
$ws.Range("AN11").Value = 5
$ws.Range("F11").Value = 1.48
$ws.Range("G11").Value = 1.68
$ws.Range("H11").Value = 4.7
$ws.Range("I11").Value = 7.6
$ws.Range("J11").Value = 3.55
$ws.Range("K11").Value = 7.2
$ws.Range("L11").Value = 1.17
$ws.Range("O11").Value = 1.11
$ws.Range("P11").Value = 1.98
$ws.Range("Q11").Value = 1.28
$ws.Range("R11").Value = 1.98
$ws.Range("S11").Value = 1.83
$ws.Range("T11").Value = 1.5
$ws.Range("U11").Value = 2.62
$ws.Range("V11").Value = 1.15
$ws.Range("W11").Value = 2.44

# Row 15
$ws.Range("F15").Value = 1.78
$ws.Range("G15").Value = 2.3
$ws.Range("I15").Value = 8.199999999999999
$ws.Range("K15").Value = 7.2
$ws.Range("L15").Value = 1.28
$ws.Range("N15").Value = 1.88
$ws.Range("P15").Value = 1.88
$ws.Range("Q15").Value = 1.65
$ws.Range("R15").Value = 1.33
$ws.Range("S15").Value = 2.62
$ws.Range("V15").Value = 1.14
$ws.Range("W15").Value = 1.76

# Row 16
$ws.Range("AC16").Value = 8
$ws.Range("AF16").Value = 14.5
$ws.Range("AG16").Value = 13
$ws.Range("AH16").Value = 25
$ws.Range("AN16").Value = 26
$ws.Range("Q16").Value = 1.02
$ws.Range("Y16").Value = 15.5
$ws.Range("Z16").Value = 38

# Row 21
$ws.Range("G21").Value = 2
$ws.Range("W21").Value = 2

# Row 23
$ws.Range("V23").Value = 1.59

# Row 24
$ws.Range("H24").Value = 3.9
$ws.Range("P24").Value = 2.28
$ws.Range("Q24").Value = 1.76

# Row 25
$ws.Range("K25").Value = 3.55

# Row 30
$ws.Range("I30").Value = 8.6
$ws.Range("U30").Value = 1.72

# Row 35
$ws.Range("G35").Value = 1.77
$ws.Range("W35").Value = 2.28

# Row 36
$ws.Range("N36").Value = 2.48

# Row 37
$ws.Range("AA37").Value = 85
$ws.Range("AB37").Value = 9
$ws.Range("AD37").Value = 16
$ws.Range("AE37").Value = 60
$ws.Range("AG37").Value = 15.5
$ws.Range("AI37").Value = 100
$ws.Range("AJ37").Value = 55
$ws.Range("AK37").Value = 48
$ws.Range("AL37").Value = 85
$ws.Range("AM37").Value = 250
$ws.Range("AN37").Value = 980
$ws.Range("AO37").Value = 90
$ws.Range("F37").Value = 2.6
$ws.Range("G37").Value = 2.72
$ws.Range("H37").Value = 3.3
$ws.Range("I37").Value = 3.55
$ws.Range("J37").Value = 2.96
$ws.Range("M37").Value = 1.13
$ws.Range("T37").Value = 2.2
$ws.Range("V37").Value = 1.4
$ws.Range("W37").Value = 1.59
$ws.Range("X37").Value = 9.800000000000001
$ws.Range("Y37").Value = 9.6
$ws.Range("Z37").Value = 26

# Row 44
$ws.Range("G44").Value = 4.2
$ws.Range("W44").Value = 1.31

# Row 45
$ws.Range("I45").Value = 2.42
$ws.Range("V45").Value = 1.7

# Row 49
$ws.Range("Q49").Value = 2.88

# Row 50
$ws.Range("H50").Value = 3.75
$ws.Range("W50").Value = 2.06

# Row 52
$ws.Range("V52").Value = 1.4

# Row 57
$ws.Range("AL57").Value = 40
